$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows above row 7 (shifts existing rows 7..100 down to 10..103)
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# Copy formatting from row 5 (A:s11 / B:s5 pair) onto the new rows
$ws.Range("A5:B5").Copy()
$ws.Range("A7:B9").PasteSpecial(-4122)

# Set the new cell values
$ws.Cells.Item(7,1).Value = "Proxy First Name"
$ws.Cells.Item(8,1).Value = "Proxy Middle Initial"
$ws.Cells.Item(9,1).Value = " Proxy Last Name"

# Row 7 should not have a B cell at all (matches target diff)
$ws.Cells.Item(7,2).Clear()

# Match the ht="19" row height used throughout the rest of the sheet
$ws.Rows.Item(7).RowHeight = 19
$ws.Rows.Item(8).RowHeight = 19
$ws.Rows.Item(9).RowHeight = 19

[void]$ws.Range("A7").Select()

Write-Output "DONE"
